# Corrected problem with direct link to Jupyter notebook
#
# Slide 3 ("Methodology") has a Content Placeholder shape whose last
# paragraph links to the analysis notebook. The link text pointed at the
# raw repo path; fix it to the real GitHub "blob" URL for the notebook.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shp = $s.Shapes.Item(5)   # "Content Placeholder 2"

$tr = $shp.TextFrame.TextRange
$fullText = $tr.Text

$oldUrl = "https://github.com/CYINT/lily-test/analysis.ipynb"
$newUrl = "https://github.com/CYINT/lily-test/blob/main/analysis.ipynb"

$startIdx = $fullText.IndexOf($oldUrl)
if ($startIdx -ge 0) {
    # TextRange.Characters uses 1-based character positions.
    $linkRange = $tr.Characters($startIdx + 1, $oldUrl.Length)
    $linkRange.Text = $newUrl
}
